$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 202.5
$ws.Range("I4").Value = 202.5
$ws.Range("K4").Value = 202.5
$ws.Range("M4").Value = -88.5
$ws.Range("H12").Value = 517.75
$ws.Range("I12").Value = 168.22223
$ws.Range("J12").Value = 1566.3334
$ws.Range("K12").Value = 168.22223
$ws.Range("L12").Value = 1566.3334
$ws.Range("M12").Value = 1.777770000000004
$ws.Range("N12").Value = -1906.3334
$ws.Range("H17").Value = 1126.9204
$ws.Range("J17").Value = 1126.9204
$ws.Range("L17").Value = 3380.7612
$ws.Range("N17").Value = -3716.7612
$ws.Range("H51").Value = 8348.5
$ws.Range("I51").Value = 13623
$ws.Range("J51").Value = 6238.7
$ws.Range("K51").Value = 13623
$ws.Range("L51").Value = 6238.7
$ws.Range("M51").Value = -13139
$ws.Range("N51").Value = -7206.7
$ws.Range("H76").Value = 336668060
$ws.Range("I76").Value = 336668060
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 336668060
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -336667745
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 336668060
$ws.Range("I79").Value = 336668060
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 336668060
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -336666968
$ws.Range("N79").ClearContents()
$ws.Range("H137").Value = 1933.7931
$ws.Range("I137").Value = 1195.0588
$ws.Range("K137").Value = 3585.1764
$ws.Range("M137").Value = -1035.1764
$ws.Range("H138").Value = 3355.388
$ws.Range("I138").Value = 2951.5356
$ws.Range("J138").Value = 3645.3333
$ws.Range("K138").Value = 8854.606800000001
$ws.Range("L138").Value = 10935.9999
$ws.Range("M138").Value = -3714.606800000001
$ws.Range("N138").Value = -21215.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3796.0833
$ws.Range("I61").Value = 4239.778
$ws.Range("K61").Value = 4239.778
$ws.Range("M61").Value = -4027.778
$ws.Range("H102").Value = 1408.1613
$ws.Range("I102").Value = 1288.4333
$ws.Range("K102").Value = 1288.4333
$ws.Range("M102").Value = 333.5667000000001
$ws.Range("H122").Value = 4539.364
$ws.Range("I122").Value = 3333
$ws.Range("K122").Value = 9999
$ws.Range("M122").Value = -7549
$ws.Range("H136").Value = 3796.0833
$ws.Range("I136").Value = 4239.778
$ws.Range("K136").Value = 12719.334
$ws.Range("M136").Value = -10169.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5957.875
$ws.Range("I86").Value = 2980.6
$ws.Range("J86").Value = 50617
$ws.Range("K86").Value = 2980.6
$ws.Range("L86").Value = 50617
$ws.Range("M86").Value = -1857.6
$ws.Range("N86").Value = -52863
$ws.Range("H89").Value = 5957.875
$ws.Range("I89").Value = 2980.6
$ws.Range("J89").Value = 50617
$ws.Range("K89").Value = 14903
$ws.Range("L89").Value = 253085
$ws.Range("M89").Value = -9287
$ws.Range("N89").Value = -264317
$ws.Range("H132").Value = 91000
$ws.Range("J132").Value = 91000
$ws.Range("L132").Value = 91000
$ws.Range("N132").Value = -101120
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H139").Value = 198199.6
$ws.Range("J139").Value = 198199.6
$ws.Range("L139").Value = 198199.6
$ws.Range("N139").Value = -208479.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1861.5358
$ws.Range("I58").Value = 1286.45
$ws.Range("J58").Value = 3299.25
$ws.Range("K58").Value = 1286.45
$ws.Range("L58").Value = 3299.25
$ws.Range("M58").Value = -1083.45
$ws.Range("N58").Value = -3705.25
$ws.Range("H99").Value = 1884.375
$ws.Range("I99").Value = 1867.8572
$ws.Range("K99").Value = 1867.8572
$ws.Range("M99").Value = -369.8571999999999
$ws.Range("H105").Value = 868.1667
$ws.Range("I105").Value = 887.6
$ws.Range("K105").Value = 887.6
$ws.Range("M105").Value = 859.4
$ws.Range("H126").Value = 1884.375
$ws.Range("I126").Value = 1867.8572
$ws.Range("K126").Value = 5603.571599999999
$ws.Range("M126").Value = -3133.571599999999
$ws.Range("H136").Value = 1861.5358
$ws.Range("I136").Value = 1286.45
$ws.Range("J136").Value = 3299.25
$ws.Range("K136").Value = 3859.35
$ws.Range("L136").Value = 9897.75
$ws.Range("M136").Value = -1309.35
$ws.Range("N136").Value = -14997.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 800
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2400
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -2624
$ws.Range("H115").Value = 1500
$ws.Range("I115").Value = 1500
$ws.Range("K115").Value = 4500
$ws.Range("M115").Value = -3325
$ws.Range("H132").Value = 2907.1
$ws.Range("I132").Value = 2824.5715
$ws.Range("J132").Value = 3099.6667
$ws.Range("K132").Value = 25421.1435
$ws.Range("L132").Value = 27897.0003
$ws.Range("M132").Value = -22891.1435
$ws.Range("N132").Value = -32957.0003
$ws.Range("H133").Value = 30000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 30000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 90000
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -100120
$ws.Range("H134").Value = 19500
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H135").Value = 800
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 800
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 7200
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -12270
$ws.Range("H136").Value = 7443
$ws.Range("I136").Value = 7443
$ws.Range("K136").Value = 22329
$ws.Range("M136").Value = -17229
$ws.Range("H137").Value = 10926
$ws.Range("I137").Value = 5815
$ws.Range("J137").Value = 14333.333
$ws.Range("K137").Value = 17445
$ws.Range("L137").Value = 42999.999
$ws.Range("M137").Value = -12345
$ws.Range("N137").Value = -53199.999
$ws.Range("H138").Value = 55570456
$ws.Range("I138").Value = 111115360
$ws.Range("K138").Value = 333346080
$ws.Range("M138").Value = -333340940
$ws.Range("H139").Value = 23824208
$ws.Range("I139").Value = 33341892
$ws.Range("K139").Value = 100025676
$ws.Range("M139").Value = -100020536
$ws.Range("H140").Value = 16679489
$ws.Range("I140").Value = 250000000
$ws.Range("J140").Value = 13738.071
$ws.Range("K140").Value = 750000000
$ws.Range("L140").Value = 41214.213
$ws.Range("M140").Value = -749994820
$ws.Range("N140").Value = -51574.213
$ws.Range("H141").Value = 22940.021
$ws.Range("I141").Value = 6095.875
$ws.Range("K141").Value = 18287.625
$ws.Range("M141").Value = -13107.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 94780
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 94780
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 94780
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -104920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 173000.5
$ws.Range("I7").Value = 255750.75
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 255750.75
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -255638.75
$ws.Range("N7").Value = -7724
$ws.Range("H40").Value = 2509.2
$ws.Range("I40").Value = 2199
$ws.Range("K40").Value = 2199
$ws.Range("M40").Value = -2063
$ws.Range("H126").Value = 173000.5
$ws.Range("I126").Value = 255750.75
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 767252.25
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -764782.25
$ws.Range("N126").Value = -27440
$ws.Range("H132").Value = 2432.111
$ws.Range("I132").Value = 2084.2632
$ws.Range("K132").Value = 6252.7896
$ws.Range("M132").Value = -3722.7896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H62").Value = 4717.636
$ws.Range("I62").Value = 4066.3333
$ws.Range("K62").Value = 4066.3333
$ws.Range("M62").Value = -3442.3333
$ws.Range("H65").Value = 4717.636
$ws.Range("I65").Value = 4066.3333
$ws.Range("K65").Value = 20331.6665
$ws.Range("M65").Value = -17211.6665
$ws.Range("H81").Value = 8550
$ws.Range("I81").Value = 8550
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 17100
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -16039
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 8550
$ws.Range("I84").Value = 8550
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 85500
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -80196
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 309886.7
$ws.Range("I122").Value = 387475.56
$ws.Range("J122").Value = 21699.428
$ws.Range("K122").Value = 1162426.68
$ws.Range("L122").Value = 65098.284
$ws.Range("M122").Value = -1159976.68
$ws.Range("N122").Value = -69998.284
$ws.Range("H126").Value = 2088.3076
$ws.Range("I126").Value = 1640
$ws.Range("K126").Value = 4920
$ws.Range("M126").Value = -2450
$ws.Range("H132").Value = 2876.7368
$ws.Range("I132").Value = 2809.4424
$ws.Range("J132").Value = 3576.6
$ws.Range("K132").Value = 8428.3272
$ws.Range("L132").Value = 10729.8
$ws.Range("M132").Value = -5898.3272
$ws.Range("N132").Value = -15789.8
